# Lab3 Part3TruthTable.xlsx - "Add finished lab3 and lab4 part 1"
# Fills in the truth-table output columns (E:K) for rows 2-17 on Sheet1
# with the completed logic-gate values, then leaves the selection where
# the author left off (cell B23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns E (col 5) through K (col 11), for rows 2-17.
$data = @(
  @(0,0,0,0,0,0,1),
  @(1,0,0,1,1,1,1),
  @(0,0,1,0,0,1,0),
  @(0,0,0,0,1,1,0),
  @(1,0,0,1,1,0,0),
  @(0,1,0,0,1,0,0),
  @(0,1,0,0,0,0,0),
  @(0,0,0,1,1,1,1),
  @(0,0,0,0,0,0,0),
  @(0,0,0,1,1,0,0),
  @(0,0,0,1,0,0,0),
  @(0,0,0,0,0,0,0),
  @(0,1,1,0,0,0,1),
  @(0,0,0,0,0,0,1),
  @(0,1,1,0,0,0,0),
  @(0,1,1,1,0,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowVals = $data[$i]
    $row = 2 + $i
    for ($j = 0; $j -lt $rowVals.Length; $j++) {
        $col = 5 + $j
        $ws.Cells.Item($row, $col).Value = $rowVals[$j]
    }
}

# Move the selection to where the author finished working.
$ws.Range("B23").Select()
